$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.2209588364395597, 0.4809798798281603, 0.4105429204137862, 0.640736233105157, 0.6087220669516331, 42)
    3  = @(0.02863981460692326, 0.4569463959652514, 0.3846277506495038, 0.6201836426813463, 0.6272182231572233, 41)
    4  = @(0.1936854491052489, 0.4757558206508753, 0.4120852984366939, 0.6419387030213196, 0.6198190349974922, 40)
    5  = @(0.04938534500615049, 0.4815667771394883, 0.4068672093313941, 0.6378614342718911, 0.6442601547138083, 39)
    6  = @(0.1962010372968331, 0.5025474175351191, 0.4435532101199133, 0.6659979054921369, 0.6449851830934898, 38)
    7  = @(0.08806700416235531, 0.491124900856406, 0.4101026362846906, 0.6403925642015924, 0.6430576491634551, 37)
    8  = @(0.1614104346065955, 0.4938957888017228, 0.4458511905737079, 0.6677208927191869, 0.657108885046729, 36)
    9  = @(0.07816625056257342, 0.5207012812240953, 0.4341256565743806, 0.6588821264644994, 0.6637803783822126, 35)
    10 = @(0.1705354281393297, 0.5116377312545397, 0.4626856766095047, 0.6802100238966673, 0.6683882402035821, 34)
    11 = @(0.05325839054106861, 0.5136449649815871, 0.4286940646171506, 0.6547473288354453, 0.6626957512451073, 33)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
